$d = $word.ActiveDocument

# --- Change 1: "Egy közvetítő osztályok keresztül ..." -> "... osztályon keresztül ..."
#     (split into three runs: "...osztályo" | "n" | " keresztül...")
$r1 = $d.Content
$found1 = $r1.Find.Execute("Egy közvetítő osztályok keresztül", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found1) {
    $start1 = $r1.Start
    # "Egy közvetítő osztályo" is 22 characters; the following "k" is the single
    # character that must become "n".
    $rK = $d.Range($start1 + 22, $start1 + 23)
    $rK.Text = "n"
    # Drop a temporary bookmark around the replaced letter and remove it again;
    # this forces Word to keep the surrounding text in its own separate run
    # (matching the three-run split in the target) without leaving any
    # formatting residue behind.
    $rMark = $d.Range($start1 + 22, $start1 + 23)
    $d.Bookmarks.Add("tmp_split_mark", $rMark)
    $d.Bookmarks("tmp_split_mark").Delete()
}

# --- Change 2: "a komponense konstruktorában" -> "a komponens konstruktorában"
$d.Content.Find.Execute("a komponense konstruktorában", $true, $false, $false, $false, $false, $true, 1, $false, "a komponens konstruktorában", 2) | Out-Null

# --- Change 3: "a többi komponense metódusai" -> "a többi komponens metódusai"
$d.Content.Find.Execute("a többi komponense metódusai", $true, $false, $false, $false, $false, $true, 1, $false, "a többi komponens metódusai", 2) | Out-Null

# --- Change 4: delete the paragraph "A nyelvtan implementálása egyszerű."
$r4 = $d.Content
$found4 = $r4.Find.Execute("A nyelvtan implementálása egyszerű.", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found4) {
    $p4 = $r4.Paragraphs(1)
    $p4.Range.Delete()
}
